$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price values that look numeric (e.g. "591.40", "40.00").
# Force text format first so Excel keeps the exact original formatting
# (trailing zeros, thousand-dot-grouped strings, etc.) instead of coercing
# them to numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.079.20"
$ws.Range("E2").Value = "  +1.27%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.145.74"
$ws.Range("E3").Value = "  +1.87%  "
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "591.40"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.43"
$ws.Range("E6").Value = "  +1.13%  "
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.135.34"
$ws.Range("E8").Value = "  +1.77%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.530"
$ws.Range("E9").Value = "  +0.74%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.162"
$ws.Range("E10").Value = "  +2.83%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.91"
$ws.Range("E11").Value = "  +5.05%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.456"
$ws.Range("E12").Value = "  +0.43%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000248"
$ws.Range("E13").Value = "  +1.31%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.25"
$ws.Range("E14").Value = "  -0.39%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.675.25"
$ws.Range("E15").Value = "  +2.02%  "
$ws.Range("E16").Value = "  -0.18%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.25"
$ws.Range("E17").Value = "  +2.16%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "63.901.00"
$ws.Range("E18").Value = "  +1.10%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.150.44"
$ws.Range("E19").Value = "  +1.87%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "467.75"
$ws.Range("E20").Value = "  +1.93%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.33"
$ws.Range("E21").Value = "  +1.17%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.731"
$ws.Range("E22").Value = "  +1.14%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.56"
$ws.Range("E23").Value = "  +2.07%  "
$ws.Range("E24").Value = "  +12.10%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.09"
$ws.Range("E25").Value = "  +1.14%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "81.01"
$ws.Range("E26").Value = "  +0.07%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.78"
$ws.Range("E28").Value = "  +9.99%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.71"
$ws.Range("E29").Value = "  +2.00%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.34"
$ws.Range("E30").Value = "  +8.50%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.22"
$ws.Range("E31").Value = "  +1.18%  "
$ws.Range("E32").Value = "  +0.12%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "27.63"
$ws.Range("E33").Value = "  +3.70%  "
$ws.Range("E34").Value = "  +4.84%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0₃0866"
$ws.Range("E35").Value = "  +2.63%  "
$ws.Range("E36").Value = "  +3.33%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.16"
$ws.Range("E37").Value = "  +2.85%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.28"
$ws.Range("E38").Value = "  -0.95%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.24"
$ws.Range("E39").Value = "  -2.28%  "
$ws.Range("B40").Value = "Cosmos"
$ws.Range("C40").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "9.40"
$ws.Range("E40").Value = "  +7.56%  "
$ws.Range("B41").Value = "Bittensor"
$ws.Range("C41").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "459.25"
$ws.Range("E41").Value = "  +5.77%  "
$ws.Range("E42").Value = "  +2.29%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.292"
$ws.Range("E43").Value = "  +9.63%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0373"
$ws.Range("E44").Value = "  +1.38%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.902.35"
$ws.Range("E45").Value = "  +1.52%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "40.00"
$ws.Range("E46").Value = "  +11.21%  "
$ws.Range("E47").Value = "  -0.40%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "133.02"
$ws.Range("E48").Value = "  +7.42%  "
$ws.Range("E49").Value = "  +0.01%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.110"
$ws.Range("E50").Value = "  +0.54%  "
$ws.Range("E51").Value = "  +4.19%  "
